# The two observation rows (row 2: Spillkråka / Dryocopus martius, and
# row 3: Mindre hackspett / Dryobates minor) have had their record-specific
# fields swapped between rows 2 and 3. Columns C, D, H, K, L, N and P..AY
# are identical between the two rows already, so only A, B, E, F, G, I, M
# need to be exchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold plain numbers/text and can be swapped via Value2
# directly (Value2 won't mis-coerce these).
$plainCols = @("A", "B", "E", "F", "G", "M")
foreach ($col in $plainCols) {
    $addr2 = "$col" + "2"
    $addr3 = "$col" + "3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $v3
    $ws.Range($addr3).Value2 = $v2
}

# Column I ("Antal") holds small digit-only strings ("1"/"2") that must
# stay text (as in the source file) rather than become numeric. Assigning
# a leading apostrophe forces Excel to keep/store it as text; resetting
# the style back to Normal afterwards drops the transient "number stored
# as text" flag Excel attaches to the cell so no visible formatting
# artifact is left behind.
$i2 = $ws.Range("I2").Value2
$i3 = $ws.Range("I3").Value2
$ws.Range("I2").Value2 = "'$i3"
$ws.Range("I3").Value2 = "'$i2"
$ws.Range("I2").Style = "Normal"
$ws.Range("I3").Style = "Normal"
